# Adds two new columns, "I0" (I) and "IF" (J), to the stats sheet.
# Header cells (I1/J1) get the same bold/centered/bordered style as the
# other header cells (B1:H1). Data cells (I2:J69) are plain values, like
# the existing data columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$headerRng = $ws.Range("I1:J1")
$headerRng.Font.Bold = $true
$headerRng.HorizontalAlignment = -4108   # xlCenter
$headerRng.VerticalAlignment = -4160     # xlTop
$headerRng.Borders.LineStyle = 1         # xlContinuous

# --- Data rows (r=2..69) ---
$I_VALS = @(4,6,7,4,8,8,7,6,5,11,9,8,6,3,6,7,5,7,4,6,3,7,4,3,8,8,6,6,8,7,8,9,6,8,9,8,7,8,8,9,7,9,6,3,7,6,7,9,9,7,8,9,7,9,9,5,8,9,9,7,6,4,9,2,5,5,3,2)
$J_VALS = @(4,6,7,4,8,8,7,6,5,11,9,8,6,4,6,7,5,7,4,6,4,7,4,3,8,8,7,6,8,7,8,9,7,8,9,8,8,8,9,9,7,9,7,3,7,6,7,9,9,8,8,9,8,9,9,6,8,9,9,8,6,4,9,2,5,5,3,2)

for ($k = 0; $k -lt $I_VALS.Length; $k++) {
    $r = 2 + $k
    $ws.Cells.Item($r, 9).Value = $I_VALS[$k]
    $ws.Cells.Item($r, 10).Value = $J_VALS[$k]
}
